# Updated cryptos list on Sat Oct 21 13:54:43 UTC 2023 with GitHub Actions
#
# Column D ("Price") values are stored as plain text in this sheet (they use
# thousands separators like "29.709.71" and Excel would otherwise try to
# coerce plain numeric-looking strings into real numbers). A leading
# apostrophe forces Excel to keep them as literal text, matching the
# original inlineStr/text storage, without altering the cell's number
# format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'29.709.71"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.607.49"
$ws.Range("E3").Value = "  +0.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'212.38"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.518"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'28.84"
$ws.Range("E8").Value = "  +7.44%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.74%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.51%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "'1.835.51"
$ws.Range("E12").Value = "  +0.27%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.618.59"
$ws.Range("E13").Value = "  +0.92%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "'0.563"
$ws.Range("E14").Value = "  +5.13%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "'29.706.71"
$ws.Range("E15").Value = "  +0.64%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +1.87%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'64.52"
$ws.Range("E17").Value = "  +1.55%  "

# Row 18/19 swap: BitcoinCash <-> Chainlink
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'8.32"
$ws.Range("E18").Value = "  +9.44%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'241.56"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.0₃0703"
$ws.Range("E20").Value = "  +1.40%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.12%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.05"
$ws.Range("E22").Value = "  +0.45%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "'9.51"
$ws.Range("E23").Value = "  +3.47%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +1.19%  "

# Row 25
$ws.Range("D25").Value = "'156.92"
$ws.Range("E25").Value = "  +1.61%  "

# Row 26
$ws.Range("D26").Value = "'15.57"
$ws.Range("E26").Value = "  +1.77%  "

# Row 27
$ws.Range("E27").Value = "  +1.09%  "

# Row 28
$ws.Range("D28").Value = "'6.51"
$ws.Range("E28").Value = "  +2.16%  "

# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.09%  "

# Row 30
$ws.Range("E30").Value = "  +1.80%  "

# Row 31
$ws.Range("E31").Value = "  +0.35%  "

# Row 32
$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  +0.55%  "

# Row 33
$ws.Range("E33").Value = "  +2.32%  "

# Row 34
$ws.Range("D34").Value = "'1.426.06"
$ws.Range("E34").Value = "  -0.14%  "

# Row 35
$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  +4.80%  "

# Row 36
$ws.Range("E36").Value = "  +0.74%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "'2.86"
$ws.Range("E37").Value = "  +1.58%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  -0.03%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +2.76%  "

# Row 40 - ImmutableX
$ws.Range("D40").Value = "'0.555"
$ws.Range("E40").Value = "  +3.92%  "

# Row 41/42 swap: RenderToken <-> ARBITRUM
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.825"
$ws.Range("E41").Value = "  +3.74%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.97"
$ws.Range("E42").Value = "  +1.03%  "

# Row 43/44 swap: Kaspa <-> BitcoinSV
$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").Value = "'54.58"
$ws.Range("E43").Value = "  +1.40%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0494"
$ws.Range("E44").Value = "  +4.71%  "

# Row 45 - PaxDollar
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'67.93"
$ws.Range("E46").Value = "  +3.56%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  +19.27%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  +2.94%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "'1.744.42"
$ws.Range("E49").Value = "  +0.16%  "

# Row 50 - Quant
$ws.Range("D50").Value = "'87.07"
$ws.Range("E50").Value = "  +0.68%  "

# Row 51 - mCoin
$ws.Range("E51").Value = "  -1.25%  "
